{"js": "// Office.js (Word JavaScript API) script.\n// Applies the \"Needs and Features\" table edits described by the diff:\n//   - Center-align the \"Priority\" and \"Planned Release\" columns for every\n//     data row (the header row is already centered).\n//   - Fill in the previously-empty \"Planned Release\" cell of the first\n//     data row with \"Before launch\".\n//   - Fill in the previously-empty \"Priority\" cell of the second data row\n//     with \"M\".\n//   - Pluralize \"app\" -> \"apps\" in the second data row's Need text so it\n//     reads \"integration with preinstalled apps like Gmail.\"\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// The \"Needs and Features\" table is the 4th table in the document (index 3).\nconst table = tables.items[3];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst PRIORITY_COL = 1;\nconst PLANNED_RELEASE_COL = 3;\n\n// Center-align the Priority and Planned Release columns for every data row\n// (row 0 is the header row, already centered).\nfor (let row = 1; row < table.rowCount; row++) {\n  for (const col of [PRIORITY_COL, PLANNED_RELEASE_COL]) {\n    const cell = table.getCell(row, col);\n    const paragraphs = cell.body.paragraphs;\n    paragraphs.load(\"items\");\n    await context.sync();\n    for (const p of paragraphs.items) {\n      p.alignment = Word.Alignment.centered;\n    }\n  }\n}\nawait context.sync();\n\n// Row 1 (\"personalised to do list.\"): Planned Release cell is empty -> \"Before launch\".\n// NOTE: use Word.InsertLocation.start (not \"replace\") so the existing (now\n// centered) paragraph is kept and just gets a run inserted into it, instead\n// of being replaced by a brand-new, left-aligned paragraph.\nconst plannedReleaseRow1 = table.getCell(1, PLANNED_RELEASE_COL);\nplannedReleaseRow1.body.insertText(\"Before launch\", Word.InsertLocation.start);\n\n// Row 2 (\"integration with preinstalled app(s) like Gmail.\"): Priority cell is empty -> \"M\".\nconst priorityRow2 = table.getCell(2, PRIORITY_COL);\npriorityRow2.body.insertText(\"M\", Word.InsertLocation.start);\n\nawait context.sync();\n\n// Row 2 Need cell: \"integration with preinstalled app like Gmail.\" ->\n// \"integration with preinstalled apps like Gmail.\"\nconst needRow2 = table.getCell(2, 0);\nconst found = needRow2.body.search(\"app\", { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\n\nif (found.items.length > 0) {\n  found.items[0].insertText(\"s\", Word.InsertLocation.end);\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies the \"Needs and Features\" table edits described by the diff:\n#   - Center-align the \"Priority\" and \"Planned Release\" columns for every\n#     data row (the header row is already centered).\n#   - Fill in the previously-empty \"Planned Release\" cell of the first\n#     data row with \"Before launch\".\n#   - Fill in the previously-empty \"Priority\" cell of the second data row\n#     with \"M\".\n#   - Pluralize \"app\" -> \"apps\" in the second data row's Need text so it\n#     reads \"integration with preinstalled apps like Gmail.\"\n\n$d = $word.ActiveDocument\n\n# The \"Needs and Features\" table is the 4th table in the document.\n$t = $d.Tables.Item(4)\n\n$wdAlignParagraphCenter = 1\n$PRIORITY_COL = 2\n$PLANNED_RELEASE_COL = 4\n\n# Center-align the Priority and Planned Release columns for every data row\n# (row 1 is the header row, already centered).\nfor ($row = 2; $row -le $t.Rows.Count; $row++) {\n  foreach ($col in @($PRIORITY_COL, $PLANNED_RELEASE_COL)) {\n    $cell = $t.Cell($row, $col)\n    $cell.Range.Paragraphs.Item(1).Alignment = $wdAlignParagraphCenter\n  }\n}\n\n# Row 2 (\"personalised to do list.\"): Planned Release cell is empty -> \"Before launch\".\n$t.Cell(2, $PLANNED_RELEASE_COL).Range.Text = \"Before launch\"\n\n# Row 3 (\"integration with preinstalled app(s) like Gmail.\"): Priority cell is empty -> \"M\".\n$t.Cell(3, $PRIORITY_COL).Range.Text = \"M\"\n\n# Row 3 Need cell: \"integration with preinstalled app like Gmail.\" ->\n# \"integration with preinstalled apps like Gmail.\"\n$needCell = $t.Cell(3, 1)\n$needCell.Range.Find.Execute(\"preinstalled app \", $false, $false, $false, $false, $false, $true, 1, $false, \"preinstalled apps \", 2) | Out-Null\n\n\"done\"\n"}
